$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Shorten column width for column A (target OOXML width="11";
# Excel's ColumnWidth->pixel->OOXML-width conversion rounds, so 10.1
# is the value that round-trips to an on-disk width of exactly 11)
$ws.Columns.Item(1).ColumnWidth = 10.1

# Update car name values in column A (rows 2-11) to shortened versions
$ws.Range("A2").Value = "Toyota"
$ws.Range("A3").Value = "Mazda"
$ws.Range("A4").Value = "Honda"
$ws.Range("A5").Value = "Land"
$ws.Range("A6").Value = "SEAT"
$ws.Range("A7").Value = "KIA"
$ws.Range("A8").Value = "Honda"
$ws.Range("A9").Value = "Hyundai"
$ws.Range("A10").Value = "ISUZU"
$ws.Range("A11").Value = "Audi"
